# Commit: "Fruta / hortaliza, semanal"
# A new weekly observation is inserted as row 27 (Berenjena, Terminal La
# Palmera de La Serena). This pushes the previously existing rows 27..133
# down to rows 28..134 (dimension grows from A1:R133 to A1:R134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 27, shifting rows 27:133 down to 28:134.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new observation's data.
$ws.Cells.Item(27, 1).Value = 8
$ws.Cells.Item(27, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44690
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = 100112001
$ws.Cells.Item(27, 7).Value = "Berenjena"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 2400
$ws.Cells.Item(27, 11).Value = 9000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 9500
$ws.Cells.Item(27, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 190
$ws.Cells.Item(27, 17).Value = 50
$ws.Cells.Item(27, 18).Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D.
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat()
